# Append newly hired giáo viên (teacher) records to the roster on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# magv, tengv, gioitinh, email, makhoa
$newTeachers = @(
    @("GV005", "Đặng Thế Hiếu",  "Nam", "dangthehieu@gmail.com", "KCB"),
    @("GV006", "Võ Văn Tuấn",    "Nam", "tuanvo@gmail.com",      "DTVT"),
    @("GV007", "Nguyễn Tuấn Anh","Nam", "tuananh@gmail.com",     "KCS"),
    @("GV008", "Nguyễn Kim Chi", "Nữ",  "kimchi@gmail.com",      "KCS"),
    @("GV009", "Bùi Công Nam",   "Nam", "congnam@gmail.com",     "CNTT-TCKGM"),
    @("GV010", "Lê Thị Giang",   "Nữ",  "legiang@gmail.com",     "CNTT-TCKGM")
)

$row = 6
foreach ($teacher in $newTeachers) {
    $ws.Range("A$row").Value = $teacher[0]
    $ws.Range("B$row").Value = $teacher[1]
    $ws.Range("C$row").Value = $teacher[2]
    $ws.Range("D$row").Value = $teacher[3]
    $ws.Range("E$row").Value = $teacher[4]
    $row++
}

$ws.Range("D14").Select()
